$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "sheet1"
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B3").Select()
